$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new rows 3-10 with the MetaDiff parameter sweep results
# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = '(0, 0)'
$ws.Range("C3").Value = 'MetaDiff'
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 32
$ws.Range("I3").Value = 50000
$ws.Range("J3").Value = 256
$ws.Range("K3").Value = 0.9097917965063573
$ws.Range("L3").Value = 0.9243668538388904
$ws.Range("M3").Value = 0.9247077910449035
$ws.Range("N3").Value = 0.9213743551510084
$ws.Range("O3").Value = 290.0301671028137
$ws.Range("P3").Value = 233.0213561058044

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = '(0, 1)'
$ws.Range("C4").Value = 'MetaDiff'
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 50000
$ws.Range("J4").Value = 256
$ws.Range("K4").Value = 0.837295951854823
$ws.Range("L4").Value = 0.9264414325164468
$ws.Range("M4").Value = 0.9267729898060869
$ws.Range("N4").Value = 0.8581380194774533
$ws.Range("O4").Value = 0.589169979095459
$ws.Range("P4").Value = 232.8850328922272

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = '(1, 0)'
$ws.Range("C5").Value = 'MetaDiff'
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 32
$ws.Range("I5").Value = 50000
$ws.Range("J5").Value = 256
$ws.Range("K5").Value = 0.9132139471504473
$ws.Range("L5").Value = 0.9307091467157134
$ws.Range("M5").Value = 0.9310223889188776
$ws.Range("N5").Value = 0.9245332691405126
$ws.Range("O5").Value = 370.4066350460052
$ws.Range("P5").Value = 239.6799750328064

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = '(1, 1)'
$ws.Range("C6").Value = 'MetaDiff'
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 32
$ws.Range("I6").Value = 50000
$ws.Range("J6").Value = 256
$ws.Range("K6").Value = 0.9140199806766977
$ws.Range("L6").Value = 0.9302733506434665
$ws.Range("M6").Value = 0.9305881187538583
$ws.Range("N6").Value = 0.9251495330680193
$ws.Range("O6").Value = 0.5930430889129639
$ws.Range("P6").Value = 202.437383890152

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = '(0, 0)'
$ws.Range("C7").Value = 'MetaDiff'
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 32
$ws.Range("I7").Value = 50000
$ws.Range("J7").Value = 256
$ws.Range("K7").Value = 0.9163044178879687
$ws.Range("L7").Value = 0.9245029929702727
$ws.Range("M7").Value = 0.9248428756753467
$ws.Range("N7").Value = 0.9270024868647453
$ws.Range("O7").Value = 280.8699870109558
$ws.Range("P7").Value = 216.8287858963013

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = '(0, 1)'
$ws.Range("C8").Value = 'MetaDiff'
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 32
$ws.Range("I8").Value = 50000
$ws.Range("J8").Value = 256
$ws.Range("K8").Value = 0.9208255071271201
$ws.Range("L8").Value = 0.9303732654229485
$ws.Range("M8").Value = 0.9306869878505435
$ws.Range("N8").Value = 0.9309828542847421
$ws.Range("O8").Value = 0.5139169692993164
$ws.Range("P8").Value = 216.0577948093414

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = '(1, 0)'
$ws.Range("C9").Value = 'MetaDiff'
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 32
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 256
$ws.Range("K9").Value = 0.9254896751661508
$ws.Range("L9").Value = 0.9323896052920266
$ws.Range("M9").Value = 0.9326939527773629
$ws.Range("N9").Value = 0.9350107156631599
$ws.Range("O9").Value = 297.2036740779877
$ws.Range("P9").Value = 221.441232919693

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = '(1, 1)'
$ws.Range("C10").Value = 'MetaDiff'
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = 10
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 32
$ws.Range("I10").Value = 50000
$ws.Range("J10").Value = 256
$ws.Range("K10").Value = 0.9282021582337634
$ws.Range("L10").Value = 0.934215142829811
$ws.Range("M10").Value = 0.9345115743203479
$ws.Range("N10").Value = 0.9374128747863907
$ws.Range("O10").Value = 0.8716640472412109
$ws.Range("P10").Value = 246.6976199150085

# Match column A styling (bold, bordered, centered) used by the header row of data (A2)
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3:A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Added rows 3-10 to Sheet1"
